# Update countries & provincias Spain
# Applies the 10-Oct-2020 00:41 COVID data refresh to the "Pais" sheet:
#  - updates the "Datos actualizados..." timestamp banner
#  - refreshes several countries' case/death figures
#  - four country pairs swap rank (adjacent rows) because the refreshed
#    totals change their sort order (list is sorted descending by
#    "Casos totales")

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Timestamp banner -------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 10 de Octubre de 2020 a las 00:41"

# --- Plain data refreshes (country stays on the same row) -------------
$ws.Range("B4").Value = 7886413
$ws.Range("C4").Value = 52650
$ws.Range("D4").Value = 5054602
$ws.Range("E4").Value = 2613266
$ws.Range("G4").Value = 807
$ws.Range("H4").Value = 218545

$ws.Range("D6").Value = 4433595
$ws.Range("E6").Value = 472654

$ws.Range("B14").Value = 688352
$ws.Range("C14").Value = 1461
$ws.Range("D14").Value = 620081
$ws.Range("E14").Value = 50724
$ws.Range("G14").Value = 139
$ws.Range("H14").Value = 17547

$ws.Range("B29").Value = 177998
$ws.Range("C29").Value = 2439
$ws.Range("D29").Value = 149420
$ws.Range("E29").Value = 18992

$ws.Range("B33").Value = 145848
$ws.Range("C33").Value = 803
$ws.Range("E33").Value = 13162
$ws.Range("G33").Value = 34
$ws.Range("H33").Value = 12175

$ws.Range("B43").Value = 104262
$ws.Range("C43").Value = 106
$ws.Range("D43").Value = 97592
$ws.Range("E43").Value = 641
$ws.Range("G43").Value = 12
$ws.Range("H43").Value = 6029

$ws.Range("B61").Value = 59992
$ws.Range("C61").Value = 151
$ws.Range("D61").Value = 51614
$ws.Range("E61").Value = 7265

$ws.Range("B74").Value = 40620
$ws.Range("C74").Value = 442
$ws.Range("D74").Value = 31876
$ws.Range("E74").Value = 7989
$ws.Range("G74").Value = 4
$ws.Range("H74").Value = 755

$ws.Range("B85").Value = 23871
$ws.Range("C85").Value = 612
$ws.Range("D85").Value = 15713
$ws.Range("E85").Value = 7271
$ws.Range("G85").Value = 7
$ws.Range("H85").Value = 887

$ws.Range("B135").Value = 4616
$ws.Range("C135").Value = 50
$ws.Range("D135").Value = 1235
$ws.Range("E135").Value = 3163
$ws.Range("G135").Value = 3
$ws.Range("H135").Value = 218

$ws.Range("B159").Value = 2051
$ws.Range("C159").Value = 1
$ws.Range("E159").Value = 129

$ws.Range("D213").Value = 12
$ws.Range("E213").Value = 4

# --- Rank swaps: the refreshed totals re-order these adjacent pairs ---

# Rows 8/9: Colombia now outranks España
$ws.Range("A8").Value = "Colombia"
$ws.Range("B8").Value = 894300
$ws.Range("C8").Value = 8121
$ws.Range("D8").Value = 780547
$ws.Range("E8").Value = 86258
$ws.Range("G8").Value = 164
$ws.Range("H8").Value = 27495

$ws.Range("A9").Value = "España"
$ws.Range("B9").Value = 890367
$ws.Range("C9").Value = 5986
$ws.Range("D9").Value = 0
$ws.Range("E9").Value = 0
$ws.Range("G9").Value = 241
$ws.Range("H9").Value = 32929

# Rows 94/95: Noruega now outranks Zambia
$ws.Range("A94").Value = "Noruega"
$ws.Range("B94").Value = 15388
$ws.Range("C94").Value = 167
$ws.Range("D94").Value = 11863
$ws.Range("E94").Value = 3250
$ws.Range("H94").Value = 275

$ws.Range("A95").Value = "Zambia"
$ws.Range("B95").Value = 15301
$ws.Range("C95").Value = 0
$ws.Range("D95").Value = 14365
$ws.Range("E95").Value = 601
$ws.Range("H95").Value = 335

# Rows 130/131: Trinidad yTobago now outranks Surinam
$ws.Range("A130").Value = "Trinidad yTobago"
$ws.Range("B130").Value = 5021
$ws.Range("C130").Value = 58
$ws.Range("D130").Value = 3150
$ws.Range("E130").Value = 1782
$ws.Range("G130").Value = 3
$ws.Range("H130").Value = 89

$ws.Range("A131").Value = "Surinam"
$ws.Range("B131").Value = 5004
$ws.Range("D131").Value = 4794
$ws.Range("E131").Value = 104
$ws.Range("H131").Value = 106

# Rows 204/205: Dominica now outranks Fiyi
$ws.Range("A204").Value = "Dominica"
$ws.Range("D204").Value = 24
$ws.Range("E204").Value = 8
$ws.Range("H204").Value = 0

$ws.Range("A205").Value = "Fiyi"
$ws.Range("B205").Value = 32
$ws.Range("D205").Value = 28
$ws.Range("E205").Value = 2
$ws.Range("H205").Value = 2
